$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 36
$ws.Range("I2").Value = 113
$ws.Range("J2").Value = 423
$ws.Range("L2").Value = 114
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 78
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 34
$ws.Range("T2").Value = 75
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 600
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 657
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 2
